$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# B11 currently holds the text "R40"; it must become the text "1".
# A bare Value = "1" would be auto-coerced to the NUMBER 1 on this
# General-formatted cell, so we use Excel's normal "force text" entry
# idiom (leading apostrophe -> quote-prefixed text) so the result stays
# a genuine string value, same as the rest of the "Rule" column (R10,
# R20, R30, ...).
$ws.Range("B11").Value = "'1"
